# Add a "Save" column (H) to the s_vals sheet, mirroring the style used
# by the other header cells (e.g. G1) and filling in the values for rows 2-7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1: same style/formatting as the other header cells (e.g. G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# Data values for H2:H7
$values = @(0, 1, 1, 0, 1, 0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
